# Insert a new record (row) into the "Hortaliza, Agrícola del Norte S.A. de Arica - Papa"
# weekly price sheet. The new row is inserted at row 98, which pushes the existing
# rows 98..165 down to 99..166 (so the sheet grows from 165 to 166 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 98, shifting everything below it down by one row.
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new price record.
$ws.Range("A98").Value = 1
$ws.Range("B98").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C98").Value = "Arica y Parinacota"
$ws.Range("D98").Value = 44981
$ws.Range("E98").Value = 15
$ws.Range("F98").Value = 100114001
$ws.Range("G98").Value = "Papa"
$ws.Range("H98").Value = "Red Lady"
$ws.Range("I98").Value = "1a (cosecha)"
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 14000
$ws.Range("L98").Value = 15000
$ws.Range("M98").Value = 14500
$ws.Range("N98").Value = "`$/saco 25 kilos"
$ws.Range("O98").Value = "Región del Maule"
$ws.Range("P98").Value = 580
$ws.Range("Q98").Value = 25
$ws.Range("R98").Value = "Hortaliza"
